# remove sy-datum/sy-uzeit from demo programs
#
# The three "per-sheet" demo worksheets (Sheet1, Sheet2, Sheet3) used to show
# the live system date/time (sy-datum / sy-uzeit) under the label
# "Current Date:". That is replaced with a fixed/static label "Date:" and a
# fixed serial value (57) instead of the live "today" date (44557). The
# navigation hyperlinks that were previously placed on these demo sheets are
# also removed.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sheet1", "Sheet2", "Sheet3")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # "Current Date:" -> "Date:"
    $ws.Range("A3").Value = "Date:"

    # Fixed integer value instead of the live sy-datum derived date serial.
    $ws.Range("A4").Value = 57

    # Drop the cross-sheet navigation hyperlink.
    $ws.Hyperlinks.Delete()
}

# Sheet4 also had a navigation hyperlink that is removed; its other content
# (including the shared-string index shift for "Date as string") is left
# untouched and handled automatically by the shared-string table rewrite.
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Hyperlinks.Delete()
